$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 54554
$ws.Range("B4").Value = 123232
$ws.Range("C7").Value = 2323232

$ws.Range("C7").Select()
